$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-24 Monday" "2025-02-25 Tuesday"

Replace-Text "23×77=" "29×60="
Replace-Text "95×54=" "91×70="
Replace-Text "65×31=" "12×86="
Replace-Text "21×35=" "13×54="
Replace-Text "41×85=" "21×90="
Replace-Text "58×51=" "99×73="
Replace-Text "26×78=" "13×47="
Replace-Text "68×60=" "35×26="
Replace-Text "95×80=" "75×41="
Replace-Text "40×85=" "58×37="
Replace-Text "54×90=" "15×79="
Replace-Text "57×28=" "52×57="
Replace-Text "51×14=" "65×60="
Replace-Text "51×88=" "20×99="
Replace-Text "61×83=" "78×88="
Replace-Text "87×66=" "95×57="
Replace-Text "95×67=" "12×15="
Replace-Text "58×97=" "55×90="
Replace-Text "82×60=" "47×43="
Replace-Text "96×81=" "33×74="
Replace-Text "27×33=" "82×89="
Replace-Text "32×28=" "43×35="
Replace-Text "71×16=" "31×65="
Replace-Text "91×57=" "78×49="
Replace-Text "99×23=" "29×64="

Write-Output "Replacements complete"
